$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at spreadsheet row 181 (pushes existing rows 181..261 down to 182..262)
$ws.Rows.Item(181).Insert()

# Populate the newly inserted row with the new weekly record
$ws.Range("A181").Value = 6
$ws.Range("B181").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C181").Value = "Metropolitana"
$ws.Range("D181").Value = 44755
$ws.Range("E181").Value = 13
$ws.Range("F181").Value = 100112026
$ws.Range("G181").Value = "Haba"
$ws.Range("H181").Value = "Sin especificar"
$ws.Range("I181").Value = "Primera"
$ws.Range("J181").Value = 270
$ws.Range("K181").Value = 16000
$ws.Range("L181").Value = 17000
$ws.Range("M181").Value = 16463
$ws.Range("N181").Value = "$/saco 25 kilos"
$ws.Range("O181").Value = "Región de Coquimbo"
$ws.Range("P181").Value = 659
$ws.Range("Q181").Value = 25
$ws.Range("R181").Value = "Hortaliza"
